$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new "Note" column before the existing "QC" column (D) ---
# This shifts the old D ("QC") and E ("QT") columns to E and F.
$ws.Columns("D").Insert()

# Match the width Excel would give a freshly authored "Note" column sized
# like the existing "Define" column (80->60 char-ish wide text column).
$ws.Columns("D").ColumnWidth = 60

# Give the new column the same per-cell formatting (borders / wrap /
# alignment) as the "Define" column it sits next to, then overwrite the
# cell contents/header for the new column.
$ws.Range("C1:C6").Copy()
$ws.Range("D1:D6").PasteSpecial(-4122)

$ws.Range("D1").Value = "Note"
$ws.Range("D2").Value = "Since the meeting takes place on Apple's property, the company has every right to make the rules.`n由于会议地点在苹果公司，公司完全有权利制定规则。"
$ws.Range("D3").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("D6").Value = ""

# --- Keep the AutoFilter / filter-database range in sync with the table ---
# Re-establish the autofilter on the new, wider header row (A1:F1) instead of
# toggling the existing one off.
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:F1").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=fruit!`$A`$1:`$F`$1"
    }
}
